$d = $word.ActiveDocument

$old = " as well as the start of the delivery date and time and end of delivery date and time."
$new = " as well as the start of the delivery date and time and end of delivery date and time, given the availability of days to deliver."

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
